# Datenbankstruktur_Arbeitsaufteilung.xlsx - "Add files via upload" edit
#
# Restructures the "Gebuchte Flüge" (booked flights) box into a "Buchungen"
# (bookings) box with a normalized FK to a brand-new "Stati" (states) lookup
# table in Y2:Z6, and nudges the sheet's view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the "Gebuchte Flüge" box title to "Buchungen"
# ---------------------------------------------------------------------
$ws.Range("L2").Value2 = "Buchungen"

# ---------------------------------------------------------------------
# 2) Header row: GID -> BID, Status -> FK SID
# ---------------------------------------------------------------------
$ws.Range("L3").Value2 = "BID"
$ws.Range("P3").Value2 = "FK SID"

# ---------------------------------------------------------------------
# 3) Data rows: the old textual Status column (Aktiv / Schon geflogen /
#    Storniert) becomes a numeric FK into the new Stati lookup table:
#      1 = Aktiv, 2 = Storniert/Geflogen-mix per row below, 3 = Geflogen
#    (values taken straight from the target workbook)
# ---------------------------------------------------------------------
$ws.Range("P4").Value2 = 1
$ws.Range("P5").Value2 = 2
$ws.Range("P6").Value2 = 3
$ws.Range("P7").Value2 = 1
$ws.Range("P8").Value2 = 2
$ws.Range("P4:P8").NumberFormat = "General"
$ws.Range("P4:P8").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 4) Brand-new "Stati" lookup box in Y2:Z6 (mirrors the look of the other
#    small boxes on the sheet: 16pt bold title, bold header row, medium
#    box border all the way around, centered content).
# ---------------------------------------------------------------------
$ws.Range("Y2").Value2 = "Stati"
$ws.Range("Y3").Value2 = "SID"
$ws.Range("Z3").Value2 = "Bez"
$ws.Range("Y4").Value2 = 1
$ws.Range("Z4").Value2 = "Aktiv "
$ws.Range("Y5").Value2 = 2
$ws.Range("Z5").Value2 = "Storniert"
$ws.Range("Y6").Value2 = 3
$ws.Range("Z6").Value2 = "Geflogen"

# Title cell Y2 (+ spacer Z2) — same look as B2 / F2 / L2 / R2
$titleRange = $ws.Range("Y2:Z2")
$titleRange.Font.Bold = $true
$titleRange.Font.Size = 16
$titleRange.HorizontalAlignment = -4108

# Header row Y3:Z3 — bold, centered
$headerRange = $ws.Range("Y3:Z3")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

# Data cells centered
$ws.Range("Y4:Z6").HorizontalAlignment = -4108

# Medium box border around the whole Y2:Z6 block
$outerBox = $ws.Range("Y2:Z6")
$outerBox.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$outerBox.Borders.Item(7).Weight = -4138  # xlMedium
$outerBox.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$outerBox.Borders.Item(8).Weight = -4138
$outerBox.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$outerBox.Borders.Item(9).Weight = -4138
$outerBox.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$outerBox.Borders.Item(10).Weight = -4138

# Vertical divider between the SID and Bez columns
$ws.Range("Y2:Y6").Borders.Item(10).LineStyle = 1
$ws.Range("Y2:Y6").Borders.Item(10).Weight = -4138

# ---------------------------------------------------------------------
# 5) View changes: selection moved, first visible column scrolled to E
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("T14").Select()
try {
    $excel.ActiveWindow.ScrollColumn = 5
} catch {
}

Write-Host "edit applied"
